# Living_Wage_Template.xlsx edit
#
# On the "LW_TW_Exhibit" sheet, insert 5 blank rows above the old row 11
# (the "Copy the values from the Area Table2..." banner row), pushing that
# banner row and the data table below it down from rows 11-18 to rows
# 16-23. Formulas in B4:C9 (which reference the data table) are
# auto-updated by the insert. The newly inserted blank rows 10-14 pick up
# the number-format styling from the row directly above them (row 9) in
# columns B:C, matching how Excel carries formatting into rows inserted
# immediately below a copied range. The chart anchored over this area is
# stretched down by the same 5 rows, and the final selection is left on
# the newly inserted row 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LW_TW_Exhibit")
$ws.Activate()

# Remember the chart's height before the insert shifts the rows it spans.
$co = $ws.ChartObjects(1)
$origChartHeight = $co.Height

# Insert 5 new rows immediately above the old row 11.
$ws.Rows("11:15").Insert()

# Carry the B9:C9 number formatting down into the newly inserted blank
# rows (B10:C14) so they match the look of the table above them.
$ws.Range("B9:C9").Copy()
$ws.Range("B10:C14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# The chart overlapping this block was anchored through old row 12; grow
# it by the 5 inserted rows (5 * 14pt row height) so it again reaches the
# same relative row/offset (now row 17) it did before the insert.
$co.Height = $origChartHeight + (5 * 14)

# Leave the selection on the newly inserted row 10, matching the editor's
# final cursor position.
$ws.Rows("10").EntireRow.Select()
